$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM values for B2:K21 (supplemental-figures update)
$data = @(
    @(-19.29293493902804, 2.456789256591561, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.352250710548291, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, 2.165544269351435, 2.851408698042349, -19.29293493902804, 2.558141906109487, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.152084454873699, -19.29293493902804),
    @(-19.29293493902804, 1.014991812038157, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.329346163096866, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(2.968762025919257, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, 2.919269606518615, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(3.605414945837459, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 1.650688062106362, -19.29293493902804, 2.19196513197457),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, 1.974025986002222, -19.29293493902804, 2.513352312867624, -19.29293493902804, -19.29293493902804, -19.29293493902804, 1.501795828822074),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, 1.658336256302025, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.203275679381774, 1.589789447089916),
    @(-19.29293493902804, -19.29293493902804, 1.706231736386457, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.067280420652906),
    @(-19.29293493902804, -19.29293493902804, -0.1917979182916582, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.379060854484933, -19.29293493902804),
    @(-19.29293493902804, 0.7635718036560544, 0.01403678352245061, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 1.119602863886341, 1.573581083268921, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 1.296128215489636, 1.470738604808477, -19.29293493902804),
    @(-19.29293493902804, -19.29293493902804, 1.742380342091904, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804, 2.024275050928313, -19.29293493902804, -19.29293493902804),
    @(-19.29293493902804, 1.583013096792767, 2.100539839365136, -19.29293493902804, 3.818695697572308, -19.29293493902804, 4.321925961123555, 1.612355835707937, -19.29293493902804, 2.432355172326707),
    @(-19.29293493902804, 1.726463717984669, -19.29293493902804, 2.419590801981661, -19.29293493902804, 3.211751444715637, -19.29293493902804, -19.29293493902804, -19.29293493902804, -19.29293493902804),
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
